$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Match the existing header formatting exactly by copying it from an existing header cell
$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Updated metric values
$ws.Range("B2").Value = 0.04959485582924521
$ws.Range("C2").Value = 0.998541499901305
$ws.Range("D2").Value = 0.1646540170263204

# Updated model description (F2), with embedded newline
$ws.Range("F2").Value = "Pipeline(steps=[('model',`n                 AdaBoostRegressor(learning_rate=0.5, n_estimators=150))])"

# New numeric cells
$ws.Range("G2").Value = 0.1260932844166139
$ws.Range("H2").Value = 0.991

# The embedded newline in F2 makes Excel auto-expand the row; auto-fit it back
# so no explicit/custom row height is persisted (matches the original sheet).
$ws.Rows(2).AutoFit()
